$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 2 (Sunbury), shifting it down to row 4
$ws.Range("A2:D3").EntireRow.Insert()
$ws.Range("A2:D3").ClearFormats()

# Row 2: Point Cook - Coffeeologist Cafe (8/2/2021)
$ws.Range("A2").Value = "Point Cook"
$ws.Range("B2").Value = "The Coffeeologist Cafe, 70/300 Point Cook Rd , Point Cook VIC 3030"
$ws.Range("C2").Value = "11:00am - 11:40am 8/2/2021"
$ws.Range("D2").Value = "Case attended venue"

# Row 3: Point Cook - Coffeeologist Cafe (10/2/2021)
$ws.Range("A3").Value = "Point Cook"
$ws.Range("B3").Value = "The Coffeeologist Cafe, 70/300 Point Cook Rd, Point Cook VIC 3030"
$ws.Range("C3").Value = "11:30am - 12:10pm 10/2/2021"
$ws.Range("D3").Value = "Case attended venue"
